# Sharon Choe resume — rewrite several "Front End Developer" bullet
# points and relocate the stray "_GoBack" bookmark from the Skills line
# to right after "Professional Experience".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Relocate the _GoBack bookmark.
#
#    A range collapsed exactly at a paragraph-final offset gets widened
#    by Bookmarks.Add to span through the following paragraph boundary,
#    so we briefly insert a marker run right after "Professional
#    Experience", anchor the (now clearly mid-paragraph) bookmark just
#    before that marker, then delete the marker again. Adding a
#    bookmark named "_GoBack" moves the existing one rather than
#    duplicating it, so the old location is cleared automatically.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Professional Experience", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter("ZZZBOOKMARKMARKERZZZ")

$rng2 = $d.Content
$rng2.Find.Execute("ZZZBOOKMARKMARKERZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rng2)

$d.Content.Find.Execute("ZZZBOOKMARKMARKERZZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# ---------------------------------------------------------------------
# 2) "Developed and designed personal webpage (Angular.JS) in
#    collaboration with back end developer." -> SharonChoe.com blurb.
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Developed and designed personal webpage (Angular.JS) in collaboration with back end developer.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "SharonChoe.com is a personal portfolio and blog built in Rails 4 and Angular JS for the front end.", `
    2)

# ---------------------------------------------------------------------
# 3) API bullet: swap "Congressman" for "GovTrack" (renaming the
#    flagged word in place keeps its spell-check proofErr wrapper) and
#    swap the trailing "WeatherApp search." wording for "Yahoo! Weather."
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    "using API calls in Congressman, Wikipedia, and ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "using API calls to ", `
    2)

$d.Content.Find.Execute("WeatherApp", $true, $false, $false, $false, $false, $true, 1, $false, "GovTrack", 2)

$d.Content.Find.Execute( `
    " search.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    " API, Wikipedia, and Yahoo! Weather.", `
    2)

# ---------------------------------------------------------------------
# 4) "Built interactive user interfaces in Simon Says and Tic Tac Toe
#    games." -> add "and designed", swap trailing wording for "(HTML/CSS)."
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    "Built interactive user interfaces", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Built and designed interactive user interfaces", `
    2)

$d.Content.Find.Execute( `
    "in Simon Says and Tic Tac Toe games.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "in Simon Says and Tic Tac Toe (HTML/CSS).", `
    2)

# ---------------------------------------------------------------------
# 5) Pomodoro Clock bullet rewrite.
# ---------------------------------------------------------------------
$d.Content.Find.Execute( `
    " called Pomodoro Clock which is a tool to improve time management. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "/time management app called Pomodoro Clock in Angular and HTML/CSS.", `
    2)

Write-Output "done"
